$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3664.4546
$ws.Range("I74").Value = 3530.9
$ws.Range("K74").Value = 3530.9
$ws.Range("M74").Value = -2594.9
$ws.Range("H77").Value = 3664.4546
$ws.Range("I77").Value = 3530.9
$ws.Range("K77").Value = 17654.5
$ws.Range("M77").Value = -12974.5
$ws.Range("H131").Value = 1536.0834
$ws.Range("I131").Value = 468.92856
$ws.Range("J131").Value = 5271.125
$ws.Range("K131").Value = 1406.78568
$ws.Range("L131").Value = 15813.375
$ws.Range("M131").Value = 3633.21432
$ws.Range("N131").Value = -25893.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 2690
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 2690
$ws.Range("M3").Value = -1885
$ws.Range("N3").Value = -2920
$ws.Range("H32").Value = 12662387
$ws.Range("I32").Value = 3833.8933
$ws.Range("K32").Value = 3833.8933
$ws.Range("M32").Value = -3546.8933
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H132").Value = 1471902.1
$ws.Range("I132").Value = 989.0606
$ws.Range("J132").Value = 8406207
$ws.Range("K132").Value = 2967.1818
$ws.Range("L132").Value = 25218621
$ws.Range("M132").Value = -437.1818000000003
$ws.Range("N132").Value = -25223681

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 3766.6667
$ws.Range("I7").Value = 650
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 650
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -537
$ws.Range("N7").Value = -10226
$ws.Range("H105").Value = 66667812
$ws.Range("I105").Value = 1247
$ws.Range("J105").Value = 200000940
$ws.Range("K105").Value = 1247
$ws.Range("L105").Value = 200000940
$ws.Range("M105").Value = 500
$ws.Range("N105").Value = -200004434

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 294
$ws.Range("I2").Value = 294
$ws.Range("K2").Value = 294
$ws.Range("M2").Value = -181
$ws.Range("H105").Value = 5242.449
$ws.Range("I105").Value = 5280.8696
$ws.Range("J105").Value = 4653.3335
$ws.Range("K105").Value = 5280.8696
$ws.Range("L105").Value = 4653.3335
$ws.Range("M105").Value = -3533.8696
$ws.Range("N105").Value = -8147.3335
$ws.Range("H132").Value = 10418226
$ws.Range("I132").Value = 1065.5714
$ws.Range("J132").Value = 30305532
$ws.Range("K132").Value = 3196.7142
$ws.Range("L132").Value = 90916596
$ws.Range("M132").Value = -666.7142000000003
$ws.Range("N132").Value = -90921656

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 33350970
$ws.Range("I44").Value = 100002680
$ws.Range("J44").Value = 25113.9
$ws.Range("K44").Value = 300008040
$ws.Range("L44").Value = 75341.70000000001
$ws.Range("M44").Value = -300007642
$ws.Range("N44").Value = -76137.70000000001
$ws.Range("H121").Value = 5051254.5
$ws.Range("I121").Value = 666
$ws.Range("J121").Value = 6536721.5
$ws.Range("K121").Value = 1998
$ws.Range("L121").Value = 19610164.5
$ws.Range("M121").Value = -688
$ws.Range("N121").Value = -19612784.5
$ws.Range("H122").Value = 7356931
$ws.Range("I122").Value = 31250322
$ws.Range("J122").Value = 5118.327
$ws.Range("K122").Value = 281252898
$ws.Range("L122").Value = 46064.943
$ws.Range("M122").Value = -281250448
$ws.Range("N122").Value = -50964.943
$ws.Range("H134").Value = 14288851
$ws.Range("I134").Value = 33333988
$ws.Range("J134").Value = 4999
$ws.Range("K134").Value = 100001964
$ws.Range("L134").Value = 14997
$ws.Range("M134").Value = -99996894
$ws.Range("N134").Value = -25137

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 12000
$ws.Range("J5").Value = 12000
$ws.Range("L5").Value = 12000
$ws.Range("N5").Value = -12224
$ws.Range("H102").Value = 2070.4285
$ws.Range("I102").Value = 2082.1667
$ws.Range("K102").Value = 2082.1667
$ws.Range("M102").Value = -460.1667000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 4000
$ws.Range("J43").Value = 4000
$ws.Range("L43").Value = 4000
$ws.Range("N43").Value = -4386
$ws.Range("H53").Value = 20000
$ws.Range("J53").Value = 20000
$ws.Range("L53").Value = 20000
$ws.Range("N53").Value = -21036
$ws.Range("H61").Value = 1384.04
$ws.Range("I61").Value = 1391.7
$ws.Range("J61").Value = 1378.9333
$ws.Range("K61").Value = 1391.7
$ws.Range("L61").Value = 1378.9333
$ws.Range("M61").Value = -1189.7
$ws.Range("N61").Value = -1782.9333
$ws.Range("H113").Value = 1384.04
$ws.Range("I113").Value = 1391.7
$ws.Range("J113").Value = 1378.9333
$ws.Range("K113").Value = 1391.7
$ws.Range("L113").Value = 1378.9333
$ws.Range("M113").Value = 778.3
$ws.Range("N113").Value = -5718.9333
$ws.Range("H132").Value = 20006272
$ws.Range("I132").Value = 38463664
$ws.Range("J132").Value = 10763.875
$ws.Range("K132").Value = 115390992
$ws.Range("L132").Value = 32291.625
$ws.Range("M132").Value = -115388462
$ws.Range("N132").Value = -37351.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1083.5555
$ws.Range("I81").Value = 921.6429000000001
$ws.Range("J81").Value = 1650.25
$ws.Range("K81").Value = 1843.2858
$ws.Range("L81").Value = 3300.5
$ws.Range("M81").Value = -782.2858000000001
$ws.Range("N81").Value = -5422.5
$ws.Range("H84").Value = 1083.5555
$ws.Range("I84").Value = 921.6429000000001
$ws.Range("J84").Value = 1650.25
$ws.Range("K84").Value = 9216.429
$ws.Range("L84").Value = 16502.5
$ws.Range("M84").Value = -3912.429
$ws.Range("N84").Value = -27110.5
$ws.Range("H131").Value = 102000
$ws.Range("J131").Value = 102000
$ws.Range("L131").Value = 102000
$ws.Range("N131").Value = -112080
$ws.Range("H132").Value = 21793.143
$ws.Range("I132").Value = 30167.473
$ws.Range("J132").Value = 6719.35
$ws.Range("K132").Value = 90502.41900000001
$ws.Range("L132").Value = 20158.05
$ws.Range("M132").Value = -87972.41900000001
$ws.Range("N132").Value = -25218.05
$ws.Range("H136").Value = 13890546
$ws.Range("I136").Value = 19231760
$ws.Range("J136").Value = 3390
$ws.Range("K136").Value = 57695280
$ws.Range("L136").Value = 10170
$ws.Range("M136").Value = -57692730
$ws.Range("N136").Value = -15270
